$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: update the Price (D) and Volume(1h) (E)
# columns with the latest scraped figures. Both columns hold plain
# text (e.g. European-style "67.804.19" prices, "  +0.10%  " deltas),
# so numeric-looking price strings are written with a leading
# apostrophe to keep Excel from reinterpreting them as numbers and
# dropping significant trailing zeros (e.g. "9.10" -> 9.1).

$ws.Range("D2").Value = "67.804.19"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.807.84"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'597.51"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'167.44"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").Value = "3.804.73"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "'6.30"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "'36.18"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "4.444.50"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "3.845.67"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "'18.56"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("D18").Value = "67.797.62"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'7.12"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "'461.53"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "'9.92"
$ws.Range("E22").Value = "  -5.82%  "
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "'0.0000155"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'83.63"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'12.12"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'10.02"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "3.953.61"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "'2.26"
$ws.Range("E32").Value = "  +5.42%  "
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").Value = "'29.77"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").Value = "'9.10"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "'3.43"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D44").Value = "'48.17"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "'150.39"
$ws.Range("E47").Value = "  +2.74%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "'398.12"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").Value = "'1.84"
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("D51").Value = "'26.36"
$ws.Range("E51").Value = "  +4.65%  "
